$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows of T2 samples to the toy dataset
$ws.Range("A4").Value = "CN"
$ws.Range("B4").Value = "F1"
$ws.Range("C4").Value = "alp"
$ws.Range("D4").Value = "r1"
$ws.Range("E4").Value = "T2"
$ws.Range("F4").Value = 5.3311498558694383

$ws.Range("A5").Value = "CN"
$ws.Range("B5").Value = "F1"
$ws.Range("C5").Value = "alp"
$ws.Range("D5").Value = "r2"
$ws.Range("E5").Value = "T2"
$ws.Range("F5").Value = 5.3617858961319023

$ws.Range("F2").Value = 5.0432241748072597
$ws.Range("F3").Value = 5.0433554250353101

$ws.Range("F11").Select()
